$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "-"

$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"

$ws.Range("C4").Value = "MEC-3B-Tec. Fundição"
$ws.Range("D4").Value = "MCT-1A-Tecnologia dos Materiais"
$ws.Range("E4").Value = "-"

$ws.Range("C6").Value = "MEC-3B-Tec. Fundição"
$ws.Range("D6").Value = "MCT-1A-Tecnologia dos Materiais"
$ws.Range("E6").Value = "-"

$ws.Range("C7").Value = "MEC-3B-Tec. Fundição"
$ws.Range("D7").Value = "MCT-1A-Tecnologia dos Materiais"
$ws.Range("E7").Value = "-"

$ws.Range("C8").Value = "MEC-3B-Tec. Fundição"
$ws.Range("D8").Value = "MCT-1A-Tecnologia dos Materiais"
$ws.Range("E8").Value = "-"

$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "MEC-3A-Tec. Fundição"

$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "MEC-3A-Tec. Fundição"

$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "MEC-3A-Tec. Fundição"

$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "MEC-3A-Tec. Fundição"

$wb.Save()
